$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds text-formatted values (commas/periods used as
# locale separators mean these are never true numbers in the source data).
# Force the Text number format first so Excel does not silently reinterpret
# numeric-looking strings (e.g. "311.74", "0.9998") as real numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.029.81'
$ws.Range("D3").Value = '1.827.04'
$ws.Range("D5").Value = '311.74'
$ws.Range("D7").Value = '0.4359'
$ws.Range("D8").Value = '0.3680'
$ws.Range("D9").Value = '0.07272'
$ws.Range("D10").Value = '0.8460'
$ws.Range("D11").Value = '20.68'
$ws.Range("D12").Value = '1.829.30'
$ws.Range("D14").Value = '0.07063'
$ws.Range("D16").Value = '89.65'
$ws.Range("D18").Value = '0.000008778'
$ws.Range("D21").Value = '27.101.46'
$ws.Range("D22").Value = '5.151'
$ws.Range("D23").Value = '10.88'
$ws.Range("D24").Value = '2.054.22'
$ws.Range("D25").Value = '1.987'
$ws.Range("D27").Value = '2.215'
$ws.Range("D29").Value = '5.236'
$ws.Range("D30").Value = '117.10'
$ws.Range("D31").Value = '0.08763'
$ws.Range("D35").Value = '4.436'
$ws.Range("D38").Value = '0.01949'
$ws.Range("D39").Value = '0.05240'
$ws.Range("D40").Value = '7.237'
$ws.Range("D43").Value = '0.5153'
$ws.Range("D45").Value = '10.66'
$ws.Range("D46").Value = '0.4775'
$ws.Range("D47").Value = '106.03'
$ws.Range("D48").Value = '1.941'
$ws.Range("D49").Value = '0.9998'
$ws.Range("D50").Value = '1.661'
$ws.Range("D51").Value = '0.06332'

# Column E ("Volume(1h)") percentage strings (kept with their original
# leading/trailing double-space padding).
$ws.Range("E2").Value = '  -1.48%  '
$ws.Range("E3").Value = '  -0.23%  '
$ws.Range("E4").Value = '  -0.29%  '
$ws.Range("E5").Value = '  -1.17%  '
$ws.Range("E6").Value = '  -0.23%  '
$ws.Range("E7").Value = '  +1.27%  '
$ws.Range("E8").Value = '  -0.77%  '
$ws.Range("E9").Value = '  +0.10%  '
$ws.Range("E10").Value = '  -2.56%  '
$ws.Range("E11").Value = '  -2.48%  '
$ws.Range("E12").Value = '  -0.14%  '
$ws.Range("E13").Value = '  -0.46%  '
$ws.Range("E14").Value = '  -0.10%  '
$ws.Range("E15").Value = '  -1.30%  '
$ws.Range("E16").Value = '  +1.89%  '
$ws.Range("E17").Value = '  -0.33%  '
$ws.Range("E19").Value = '  -0.24%  '
$ws.Range("E20").Value = '  -2.32%  '
$ws.Range("E21").Value = '  -1.25%  '
$ws.Range("E22").Value = '  -0.53%  '
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("E24").Value = '  -0.40%  '
$ws.Range("E25").Value = '  -1.52%  '
$ws.Range("E26").Value = '  -1.20%  '
$ws.Range("E27").Value = '  +3.21%  '
$ws.Range("E28").Value = '  -0.83%  '
$ws.Range("E29").Value = '  -1.40%  '
$ws.Range("E30").Value = '  -0.42%  '
$ws.Range("E31").Value = '  -0.97%  '
$ws.Range("E32").Value = '  -2.64%  '
$ws.Range("E33").Value = '  -3.65%  '
$ws.Range("E34").Value = '  +0.24%  '
$ws.Range("E35").Value = '  -1.64%  '
$ws.Range("E36").Value = '  -0.29%  '
$ws.Range("E37").Value = '  -2.41%  '
$ws.Range("E39").Value = '  -0.99%  '
$ws.Range("E40").Value = '  +0.84%  '
$ws.Range("E41").Value = '  -0.54%  '
$ws.Range("E42").Value = '  +1.17%  '
$ws.Range("E43").Value = '  +1.27%  '
$ws.Range("E44").Value = '  -1.42%  '
$ws.Range("E45").Value = '  -0.07%  '
$ws.Range("E46").Value = '  +0.68%  '
$ws.Range("E47").Value = '  -0.56%  '
$ws.Range("E48").Value = '  +5.82%  '
$ws.Range("E49").Value = '  -0.31%  '
$ws.Range("E50").Value = '  -0.83%  '
$ws.Range("E51").Value = '  -1.42%  '
